{"js": "// FPA Application Form - File Number fix\n// The \"File No.\" field reads \".../FPA/<year>/SDFPB/IDK\" - the trailing\n// \"IDK\" placeholder is replaced by a bookmarked \"Unit\" placeholder\n// (named \"District\"), matching the existing \"year1\"/\"year2\" placeholder\n// pattern used elsewhere in the form.\n\n// 1. Find the run of text that ends with the \"IDK\" placeholder and trim\n//    it back down to just \"/SDFPB/\" (keeps the surrounding run's\n//    formatting since we're editing in place).\nconst hits = context.document.body.search(\"/SDFPB/IDK\", { matchCase: true });\nhits.load(\"text\");\nawait context.sync();\n\nconst oldRange = hits.items[0];\noldRange.insertText(\"/SDFPB/\", \"Replace\");\nawait context.sync();\n\n// 2. Insert the new \"Unit\" placeholder text right after \"/SDFPB/\" and\n//    wrap it in a bookmark named \"District\", the same way\n//    \"year1\"/\"year2\"/\"Office\" are bookmarked elsewhere in the document.\nconst stem = context.document.body.search(\"/SDFPB/\", { matchCase: true });\nstem.load(\"text\");\nawait context.sync();\n\nconst stemRange = stem.items[0];\nconst unitRange = stemRange.insertText(\"Unit\", \"After\");\nawait context.sync();\n\nunitRange.insertBookmark(\"District\");\nawait context.sync();\n", "ps1": "# FPA Application Form - File Number fix\n# The \"File No.\" field reads \".../FPA/<year>/SDFPB/IDK\" - the trailing\n# \"IDK\" placeholder is replaced by a bookmarked \"Unit\" placeholder\n# (named \"District\"), matching the existing \"year1\"/\"year2\" placeholder\n# pattern used elsewhere in the form.\n\n$d = $word.ActiveDocument\n\n# 1. Split \"/SDFPB/IDK\" into \"/SDFPB/\" + \"Unit\" by replacing the trailing\n#    \"IDK\" text in place (keeps the surrounding run's formatting).\n$r = $d.Content\n$r.Find.Execute(\"/SDFPB/IDK\")\n$r.Text = \"/SDFPB/Unit\"\n\n# 2. Wrap the new \"Unit\" placeholder in a bookmark named \"District\", the\n#    same way \"year1\"/\"year2\"/\"Office\" are bookmarked elsewhere in the\n#    document. Word automatically renumbers any existing bookmark ids\n#    that collide with the newly inserted one.\n$r2 = $d.Content\n$r2.Find.Execute(\"Unit\")\n$d.Bookmarks.Add(\"District\", $r2)\n"}
